$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 ("time_taken"), matching the style of the other
# header cells (E1) by copying its format, then overwriting the value.
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Add the time_taken data values for rows 2 and 3 (plain/default style).
$ws.Range("F2").Value = "2021-10-05 10:50:00.212503"
$ws.Range("F3").Value = "2021-10-05 10:50:00.212514"
